$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G2").Value = 115.2213693333333
$ws.Range("H2").Value = 345.664108
$ws.Range("I2").Value = 0.2787408744545015
$ws.Range("J2").Value = 0.2787408744545015
$ws.Range("M2").Value = 8.461686666666667
$ws.Range("N2").Value = 25.38506
$ws.Range("O2").Value = 0.1873739652872041
$ws.Range("P2").Value = 0.1873739652872041
$ws.Range("Q2").Value = 974.9671246029422
$ws.Range("R2").Value = 8774.704121426479
$ws.Range("S2").Value = 0.05222878293416269
$ws.Range("T2").Value = 0.05222878293416268

$ws.Range("G3").Value = 115.2213693333333
$ws.Range("H3").Value = 345.664108
$ws.Range("I3").Value = 0.2787408744545015
$ws.Range("J3").Value = 0.2787408744545015
$ws.Range("O3").Value = 0.02691675086296081
$ws.Range("P3").Value = 0.02691675086296081
$ws.Range("Q3").Value = 140.056528943548
$ws.Range("R3").Value = 1260.508760491932
$ws.Range("S3").Value = 0.007502798673015653
$ws.Range("T3").Value = 0.007502798673015653

$ws.Range("G4").Value = 115.2213693333333
$ws.Range("H4").Value = 345.664108
$ws.Range("I4").Value = 0.2787408744545015
$ws.Range("J4").Value = 0.2787408744545015
$ws.Range("M4").Value = 16.89366666666666
$ws.Range("N4").Value = 50.681
$ws.Range("O4").Value = 0.3740901118500721
$ws.Range("P4").Value = 0.3740901118500721
$ws.Range("Q4").Value = 1946.511406394222
$ws.Range("R4").Value = 17518.602657548
$ws.Range("S4").Value = 0.1042742049018714
$ws.Range("T4").Value = 0.1042742049018714

$ws.Range("G5").Value = 115.2213693333333
$ws.Range("H5").Value = 345.664108
$ws.Range("I5").Value = 0.2787408744545015
$ws.Range("J5").Value = 0.2787408744545015
$ws.Range("M5").Value = 18.58845466666667
$ws.Range("N5").Value = 55.76536400000001
$ws.Range("O5").Value = 0.4116191719997629
$ws.Range("P5").Value = 0.411619171999763
$ws.Range("Q5").Value = 2141.787200483924
$ws.Range("R5").Value = 19276.08480435531
$ws.Range("S5").Value = 0.1147350879454518
$ws.Range("T5").Value = 0.1147350879454518

$ws.Range("I6").Value = 0.44716501655323
$ws.Range("J6").Value = 0.4471650165532299
$ws.Range("M6").Value = 8.461686666666667
$ws.Range("N6").Value = 25.38506
$ws.Range("O6").Value = 0.1873739652872041
$ws.Range("P6").Value = 0.1873739652872041
$ws.Range("Q6").Value = 1564.07341142604
$ws.Range("R6").Value = 14076.66070283436
$ws.Range("S6").Value = 0.08378708228929697
$ws.Range("T6").Value = 0.08378708228929696

$ws.Range("I7").Value = 0.44716501655323
$ws.Range("J7").Value = 0.4471650165532299
$ws.Range("O7").Value = 0.02691675086296081
$ws.Range("P7").Value = 0.02691675086296081
$ws.Range("S7").Value = 0.01203622934519504
$ws.Range("T7").Value = 0.01203622934519504

$ws.Range("I8").Value = 0.44716501655323
$ws.Range("J8").Value = 0.4471650165532299
$ws.Range("M8").Value = 16.89366666666666
$ws.Range("N8").Value = 50.681
$ws.Range("O8").Value = 0.3740901118500721
$ws.Range("P8").Value = 0.3740901118500721
$ws.Range("Q8").Value = 3122.655789054
$ws.Range("R8").Value = 28103.902101486
$ws.Range("S8").Value = 0.1672800110578371
$ws.Range("T8").Value = 0.1672800110578371

$ws.Range("I9").Value = 0.44716501655323
$ws.Range("J9").Value = 0.4471650165532299
$ws.Range("M9").Value = 18.58845466666667
$ws.Range("N9").Value = 55.76536400000001
$ws.Range("O9").Value = 0.4116191719997629
$ws.Range("P9").Value = 0.411619171999763
$ws.Range("Q9").Value = 3435.923456981976
$ws.Range("R9").Value = 30923.31111283779
$ws.Range("S9").Value = 0.1840616938609008
$ws.Range("T9").Value = 0.1840616938609008

$ws.Range("G10").Value = 60.55095666666667
$ws.Range("H10").Value = 181.65287
$ws.Range("I10").Value = 0.1464834753134679
$ws.Range("J10").Value = 0.1464834753134678
$ws.Range("M10").Value = 8.461686666666667
$ws.Range("N10").Value = 25.38506
$ws.Range("O10").Value = 0.1873739652872041
$ws.Range("P10").Value = 0.1873739652872041
$ws.Range("Q10").Value = 512.3632226802445
$ws.Range("R10").Value = 4611.2690041222
$ws.Range("S10").Value = 0.02744718961853475
$ws.Range("T10").Value = 0.02744718961853474

$ws.Range("G11").Value = 60.55095666666667
$ws.Range("H11").Value = 181.65287
$ws.Range("I11").Value = 0.1464834753134679
$ws.Range("J11").Value = 0.1464834753134678
$ws.Range("O11").Value = 0.02691675086296081
$ws.Range("P11").Value = 0.02691675086296081
$ws.Range("Q11").Value = 73.60229151947001
$ws.Range("R11").Value = 662.42062367523
$ws.Range("S11").Value = 0.003942859210553284
$ws.Range("T11").Value = 0.003942859210553284

$ws.Range("G12").Value = 60.55095666666667
$ws.Range("H12").Value = 181.65287
$ws.Range("I12").Value = 0.1464834753134679
$ws.Range("J12").Value = 0.1464834753134678
$ws.Range("M12").Value = 16.89366666666666
$ws.Range("N12").Value = 50.681
$ws.Range("O12").Value = 0.3740901118500721
$ws.Range("P12").Value = 0.3740901118500721
$ws.Range("Q12").Value = 1022.927678274444
$ws.Range("R12").Value = 9206.34910447
$ws.Range("S12").Value = 0.05479801966420247
$ws.Range("T12").Value = 0.05479801966420246

$ws.Range("G13").Value = 60.55095666666667
$ws.Range("H13").Value = 181.65287
$ws.Range("I13").Value = 0.1464834753134679
$ws.Range("J13").Value = 0.1464834753134678
$ws.Range("M13").Value = 18.58845466666667
$ws.Range("N13").Value = 55.76536400000001
$ws.Range("O13").Value = 0.4116191719997629
$ws.Range("P13").Value = 0.411619171999763
$ws.Range("Q13").Value = 1125.548713021631
$ws.Range("R13").Value = 10129.93841719468
$ws.Range("S13").Value = 0.06029540682017736
$ws.Range("T13").Value = 0.06029540682017735

$ws.Range("G14").Value = 52.74960833333333
$ws.Range("H14").Value = 158.248825
$ws.Range("I14").Value = 0.1276106336788006
$ws.Range("J14").Value = 0.1276106336788006
$ws.Range("M14").Value = 8.461686666666667
$ws.Range("N14").Value = 25.38506
$ws.Range("O14").Value = 0.1873739652872041
$ws.Range("P14").Value = 0.1873739652872041
$ws.Range("Q14").Value = 446.3506575060556
$ws.Range("R14").Value = 4017.1559175545
$ws.Range("S14").Value = 0.02391091044520971
$ws.Range("T14").Value = 0.02391091044520971

$ws.Range("G15").Value = 52.74960833333333
$ws.Range("H15").Value = 158.248825
$ws.Range("I15").Value = 0.1276106336788006
$ws.Range("J15").Value = 0.1276106336788006
$ws.Range("O15").Value = 0.02691675086296081
$ws.Range("P15").Value = 0.02691675086296081
$ws.Range("Q15").Value = 64.11941716232501
$ws.Range("R15").Value = 577.0747544609251
$ws.Range("S15").Value = 0.003434863634196832
$ws.Range("T15").Value = 0.003434863634196832

$ws.Range("G16").Value = 52.74960833333333
$ws.Range("H16").Value = 158.248825
$ws.Range("I16").Value = 0.1276106336788006
$ws.Range("J16").Value = 0.1276106336788006
$ws.Range("M16").Value = 16.89366666666666
$ws.Range("N16").Value = 50.681
$ws.Range("O16").Value = 0.3740901118500721
$ws.Range("P16").Value = 0.3740901118500721
$ws.Range("Q16").Value = 891.1342999805555
$ws.Range("R16").Value = 8020.208699825
$ws.Range("S16").Value = 0.0477378762261611
$ws.Range("T16").Value = 0.04773787622616111

$ws.Range("G17").Value = 52.74960833333333
$ws.Range("H17").Value = 158.248825
$ws.Range("I17").Value = 0.1276106336788006
$ws.Range("J17").Value = 0.1276106336788006
$ws.Range("M17").Value = 18.58845466666667
$ws.Range("N17").Value = 55.76536400000001
$ws.Range("O17").Value = 0.4116191719997629
$ws.Range("P17").Value = 0.411619171999763
$ws.Range("Q17").Value = 980.5337031885889
$ws.Range("R17").Value = 8824.803328697302
$ws.Range("S17").Value = 0.05252698337323298
$ws.Range("T17").Value = 0.05252698337323299
